# The deck currently uses the "Integral" design (theme) on its slide
# master. The author switched the presentation's design back to the
# built-in Office "Office Theme" colour scheme (Design tab > Themes >
# Office). Re-apply that by writing the 12 standard Office theme colours
# into the active theme's colour scheme, in MsoThemeColorSchemeIndex
# order: Dark1, Light1, Dark2, Light2, Accent1-6, Hyperlink,
# FollowedHyperlink.

$p = $ppt.ActivePresentation

$theme = $p.SlideMaster.Theme

$officeThemeColors = @(
    "000000",   # 1  Dark 1
    "FFFFFF",   # 2  Light 1
    "44546A",   # 3  Dark 2
    "E7E6E6",   # 4  Light 2
    "5B9BD5",   # 5  Accent 1
    "ED7D31",   # 6  Accent 2
    "A5A5A5",   # 7  Accent 3
    "FFC000",   # 8  Accent 4
    "4472C4",   # 9  Accent 5
    "70AD47",   # 10 Accent 6
    "0563C1",   # 11 Hyperlink
    "954F72"    # 12 Followed Hyperlink
)

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $hex = $officeThemeColors[$i]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)

    # PowerPoint's RGB colour values are packed as 0x00BBGGRR.
    $rgbValue = ($b * 65536) + ($g * 256) + $r

    $theme.ThemeColorScheme.Item($i + 1).RGB = $rgbValue
}
